$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the team record columns (AD:AF), appended after
# the existing "Unnamed: 28" column (AC).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, thin border) from an
# existing header cell so the new headers match the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (Wins/Losses/Ties) for every player row.
$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 96
    $ws.Cells.Item($r, 31).Value = 66
    $ws.Cells.Item($r, 32).Value = 0
}
